$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New row 11: Name/Value pair, matching the formatting pattern of the
# existing Name/Value rows (style copied from A7, which has the same
# "label" styling used throughout column A).
$ws.Range("A7").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "ExtractLimitFromURL"
$ws.Range("B11").Value = 5

# Move the active selection to the newly added cell, as in the source workbook
$ws.Range("A11").Select()
